# Add a "Kode" column header next to the existing No / Nama / No. WA headers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Kode"

# Mirror the author's final selection (cell A2) left after editing.
$ws.Range("A2").Select()
